$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1034.75
$ws.Range("J17").Value = 1034.75
$ws.Range("L17").Value = 3104.25
$ws.Range("N17").Value = -3440.25
$ws.Range("H32").Value = 3204
$ws.Range("I32").Value = 2932.6667
$ws.Range("J32").Value = 3366.8
$ws.Range("K32").Value = 2932.6667
$ws.Range("L32").Value = 3366.8
$ws.Range("M32").Value = -2606.6667
$ws.Range("N32").Value = -4018.8
$ws.Range("H40").Value = 1886.742
$ws.Range("I40").Value = 1799.6666
$ws.Range("J40").Value = 2185.2856
$ws.Range("K40").Value = 1799.6666
$ws.Range("L40").Value = 2185.2856
$ws.Range("M40").Value = -1624.6666
$ws.Range("N40").Value = -2535.2856
$ws.Range("H47").Value = 17500
$ws.Range("I47").Value = 10000
$ws.Range("J47").Value = 25000
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 25000
$ws.Range("M47").Value = -9028
$ws.Range("N47").Value = -26944
$ws.Range("H125").Value = 3556
$ws.Range("I125").Value = 3350
$ws.Range("J125").Value = 4998
$ws.Range("K125").Value = 30150
$ws.Range("L125").Value = 44982
$ws.Range("M125").Value = -27690
$ws.Range("N125").Value = -49902
$ws.Range("H126").Value = 100000
$ws.Range("J126").Value = 100000
$ws.Range("L126").Value = 100000
$ws.Range("N126").Value = -109880
$ws.Range("H137").Value = 1212.2142
$ws.Range("I137").Value = 1000.8
$ws.Range("K137").Value = 3002.4
$ws.Range("M137").Value = -452.3999999999996

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H88").Value = 2761.6365
$ws.Range("I88").Value = 2675.25
$ws.Range("J88").Value = 2811
$ws.Range("K88").Value = 2675.25
$ws.Range("L88").Value = 2811
$ws.Range("M88").Value = -2269.25
$ws.Range("N88").Value = -3623
$ws.Range("H91").Value = 2761.6365
$ws.Range("I91").Value = 2675.25
$ws.Range("J91").Value = 2811
$ws.Range("K91").Value = 2675.25
$ws.Range("L91").Value = 2811
$ws.Range("M91").Value = -1271.25
$ws.Range("N91").Value = -5619
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802
$ws.Range("H97").Value = 900.1
$ws.Range("I97").Value = 850.25
$ws.Range("K97").Value = 850.25
$ws.Range("M97").Value = -354.25
$ws.Range("H122").Value = 29725.584
$ws.Range("I122").Value = 29725.584
$ws.Range("K122").Value = 89176.75199999999
$ws.Range("M122").Value = -86726.75199999999
$ws.Range("H132").Value = 1200
$ws.Range("I132").Value = 1200
$ws.Range("K132").Value = 3600
$ws.Range("M132").Value = -1070
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H138").Value = 2449964.5
$ws.Range("J138").Value = 2449964.5
$ws.Range("L138").Value = 2449964.5
$ws.Range("N138").Value = -2460244.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 989.125
$ws.Range("I86").Value = 844.7143
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 844.7143
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = 278.2857
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 989.125
$ws.Range("I89").Value = 844.7143
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 4223.5715
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = 1392.4285
$ws.Range("N89").Value = -21232

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 993.8333
$ws.Range("I10").Value = 993.8333
$ws.Range("K10").Value = 993.8333
$ws.Range("M10").Value = -854.8333
$ws.Range("H22").Value = 1320.8572
$ws.Range("I22").Value = 1229.8
$ws.Range("J22").Value = 1548.5
$ws.Range("K22").Value = 1229.8
$ws.Range("L22").Value = 1548.5
$ws.Range("M22").Value = -879.8
$ws.Range("N22").Value = -2248.5
$ws.Range("H31").Value = 1457.5588
$ws.Range("I31").Value = 1049.6923
$ws.Range("J31").Value = 1710.0476
$ws.Range("K31").Value = 1049.6923
$ws.Range("L31").Value = 1710.0476
$ws.Range("M31").Value = -754.6922999999999
$ws.Range("N31").Value = -2300.0476
$ws.Range("H34").Value = 1457.5588
$ws.Range("I34").Value = 1049.6923
$ws.Range("J34").Value = 1710.0476
$ws.Range("K34").Value = 1049.6923
$ws.Range("L34").Value = 1710.0476
$ws.Range("M34").Value = -847.6922999999999
$ws.Range("N34").Value = -2114.0476
$ws.Range("H74").Value = 70000
$ws.Range("J74").Value = 70000
$ws.Range("L74").Value = 70000
$ws.Range("N74").Value = -71748
$ws.Range("H77").Value = 70000
$ws.Range("J77").Value = 70000
$ws.Range("L77").Value = 210000
$ws.Range("N77").Value = -218736
$ws.Range("H132").Value = 4636.273
$ws.Range("I132").Value = 4636.273
$ws.Range("K132").Value = 13908.819
$ws.Range("M132").Value = -11378.819

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 5917.143
$ws.Range("I6").Value = 8271.4
$ws.Range("J6").Value = 31.5
$ws.Range("K6").Value = 24814.2
$ws.Range("L6").Value = 94.5
$ws.Range("M6").Value = -24701.2
$ws.Range("N6").Value = -320.5
$ws.Range("H14").Value = 4999
$ws.Range("I14").Value = 4999
$ws.Range("K14").Value = 14997
$ws.Range("M14").Value = -14824
$ws.Range("H99").Value = 385
$ws.Range("J99").Value = 750
$ws.Range("L99").Value = 2250
$ws.Range("N99").Value = -6742
$ws.Range("H109").Value = 899
$ws.Range("I109").Value = 899
$ws.Range("K109").Value = 2697
$ws.Range("M109").Value = -1657
$ws.Range("H134").Value = 1264.6428
$ws.Range("I134").Value = 1238.8462
$ws.Range("K134").Value = 3716.5386
$ws.Range("M134").Value = 1353.4614

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 133.7
$ws.Range("I2").Value = 135
$ws.Range("J2").Value = 126.333336
$ws.Range("K2").Value = 135
$ws.Range("L2").Value = 126.333336
$ws.Range("M2").Value = -22
$ws.Range("N2").Value = -352.333336
$ws.Range("H97").Value = 366.66666
$ws.Range("I97").Value = 300
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 300
$ws.Range("L97").Value = 400
$ws.Range("M97").Value = 196
$ws.Range("N97").Value = -1392
$ws.Range("H122").Value = 2359.625
$ws.Range("I122").Value = 2227.4614
$ws.Range("K122").Value = 6682.3842
$ws.Range("M122").Value = -4232.3842
$ws.Range("H126").Value = 6125
$ws.Range("I126").Value = 4370
$ws.Range("K126").Value = 13110
$ws.Range("M126").Value = -10640
$ws.Range("H132").Value = 6054.2
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 9000
$ws.Range("N132").Value = -14060
$ws.Range("H135").Value = 255000
$ws.Range("J135").Value = 255000
$ws.Range("L135").Value = 255000
$ws.Range("N135").Value = -265140

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5205.778
$ws.Range("I7").Value = 3190.8
$ws.Range("K7").Value = 3190.8
$ws.Range("M7").Value = -3078.8
$ws.Range("H40").Value = 2825.9285
$ws.Range("I40").Value = 2034
$ws.Range("K40").Value = 2034
$ws.Range("M40").Value = -1898
$ws.Range("H46").Value = 1805.25
$ws.Range("J46").Value = 1736.2727
$ws.Range("L46").Value = 1736.2727
$ws.Range("N46").Value = -2112.2727
$ws.Range("H55").Value = 1560.0625
$ws.Range("I55").Value = 1559.5
$ws.Range("K55").Value = 1559.5
$ws.Range("M55").Value = -1386.5
$ws.Range("H68").Value = 2171
$ws.Range("I68").Value = 2092.5
$ws.Range("J68").Value = 2249.5
$ws.Range("K68").Value = 2092.5
$ws.Range("L68").Value = 2249.5
$ws.Range("M68").Value = -1343.5
$ws.Range("N68").Value = -3747.5
$ws.Range("H71").Value = 2171
$ws.Range("I71").Value = 2092.5
$ws.Range("J71").Value = 2249.5
$ws.Range("K71").Value = 10462.5
$ws.Range("L71").Value = 11247.5
$ws.Range("M71").Value = -6718.5
$ws.Range("N71").Value = -18735.5
$ws.Range("H76").Value = 33429
$ws.Range("J76").Value = 33429
$ws.Range("L76").Value = 33429
$ws.Range("N76").Value = -34105
$ws.Range("H79").Value = 33429
$ws.Range("J79").Value = 33429
$ws.Range("L79").Value = 33429
$ws.Range("N79").Value = -35769
$ws.Range("H100").Value = 2499.5
$ws.Range("I100").Value = 2499.5
$ws.Range("K100").Value = 2499.5
$ws.Range("M100").Value = -1958.5
$ws.Range("H126").Value = 5205.778
$ws.Range("I126").Value = 3190.8
$ws.Range("K126").Value = 9572.400000000001
$ws.Range("M126").Value = -7102.400000000001
$ws.Range("H132").Value = 2438.6086
$ws.Range("I132").Value = 2095.6667
$ws.Range("K132").Value = 6287.000100000001
$ws.Range("M132").Value = -3757.000100000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 4000
$ws.Range("I18").Value = 4000
$ws.Range("K18").Value = 4000
$ws.Range("M18").Value = -3827
$ws.Range("H31").Value = 21000
$ws.Range("J31").Value = 21000
$ws.Range("L31").Value = 21000
$ws.Range("N31").Value = -21696
$ws.Range("H122").Value = 4406.75
$ws.Range("I122").Value = 3321.6365
$ws.Range("J122").Value = 6794
$ws.Range("K122").Value = 9964.9095
$ws.Range("L122").Value = 20382
$ws.Range("M122").Value = -7514.9095
$ws.Range("N122").Value = -25282
